$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = ""
$ws.Range("C8").Value = 33
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.0"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it's ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.00"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "P. point"
$ws.Range("C9").Value = 37
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F9").Value = 472
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "17464.00"
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = ""
$ws.Range("C10").Value = 78
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.0"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F10").Value = 0
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.00"
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "P. point"
$ws.Range("C11").Value = 78
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "On board"
$ws.Range("F11").Value = 136
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "10608.00"
$ws.Range("C12").Value = 71
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.0"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F12").Value = 23
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "1633.00"
$ws.Range("C13").Value = 72
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.0"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F13").Value = 50
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "3600.00"
$ws.Range("C14").Value = 50
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.0"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F14").Value = 33
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "1650.00"
$ws.Range("C15").Value = 28
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F15").Value = 30
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "840.00"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "Each"
$ws.Range("C16").Value = 67
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.0"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F16").Value = 30
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "2010.00"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "Each"
$ws.Range("C17").Value = 32
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.0"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F17").Value = 219
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "7008.00"
$ws.Range("C18").Value = 82
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "R. mtr."
$ws.Range("C19").Value = 42
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "20 mm"
$ws.Range("F19").Value = 40
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "1680.00"
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "R. mtr."
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "25 mm"
$ws.Range("F20").Value = 56
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "4872.00"
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = ""
$ws.Range("C21").Value = 35
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.0"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F21").Value = 0
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0.00"
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "Mtr."
$ws.Range("C22").Value = 66
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "20"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2 x 4.0 sq. mm. + 1 x 2.5 sq. mm."
$ws.Range("F22").Value = 122
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "8052.00"
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = ""
$ws.Range("C23").Value = 57
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("F23").Value = 0
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0.00"
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "Mtr."
$ws.Range("C24").Value = 55
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "8 SWG G.I. ( Hot Dipped  ) Wire "
$ws.Range("F24").Value = 20
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "1100.00"
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = ""
$ws.Range("C25").Value = 64
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.0"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F25").Value = 0
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "0.00"
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = ""
$ws.Range("C26").Value = 69
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.0"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F26").Value = 0
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "0.00"
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "Each"
$ws.Range("C27").Value = 55
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F27").Value = 492
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "27060.00"
$ws.Range("C28").Value = 33
$ws.Range("C29").Value = 6
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "1122.00"
$ws.Range("C30").Value = 17
$ws.Range("C31").Value = 96
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "86400.00"
$ws.Range("C32").Value = 84
$ws.Range("C33").Value = 6
$ws.Range("C34").Value = 22
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "48048.00"
$ws.Range("C35").Value = 3
$ws.Range("C36").Value = 64
$ws.Range("C37").Value = 97
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "223147.00"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "223147.00"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "223147.00"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "223147.00"

Write-Host "Applied all cell updates"